{"js": "// Replace the two-digit multiplication answers throughout the document body.\n// Each old string is unique within the document, so a straightforward\n// search-and-replace (matching whole-word, case-sensitive) for each pair\n// is sufficient and safe.\nconst replacements = [\n  [\"19\u00d752=988\", \"33\u00d753=1749\"],\n  [\"27\u00d727=729\", \"85\u00d788=7480\"],\n  [\"17\u00d756=952\", \"79\u00d787=6873\"],\n  [\"44\u00d760=2640\", \"73\u00d718=1314\"],\n  [\"63\u00d728=1764\", \"21\u00d765=1365\"],\n  [\"91\u00d754=4914\", \"60\u00d786=5160\"],\n  [\"85\u00d732=2720\", \"35\u00d789=3115\"],\n  [\"65\u00d765=4225\", \"68\u00d747=3196\"],\n  [\"64\u00d720=1280\", \"17\u00d781=1377\"],\n  [\"32\u00d718=576\", \"34\u00d767=2278\"],\n  [\"56\u00d743=2408\", \"12\u00d765=780\"],\n  [\"69\u00d712=828\", \"75\u00d743=3225\"],\n  [\"39\u00d731=1209\", \"14\u00d757=798\"],\n  [\"13\u00d781=1053\", \"52\u00d771=3692\"],\n  [\"39\u00d750=1950\", \"72\u00d713=936\"],\n  [\"81\u00d758=4698\", \"97\u00d784=8148\"],\n  [\"36\u00d738=1368\", \"44\u00d766=2904\"],\n  [\"15\u00d760=900\", \"63\u00d785=5355\"],\n  [\"39\u00d722=858\", \"82\u00d742=3444\"],\n  [\"37\u00d757=2109\", \"39\u00d759=2301\"],\n  [\"98\u00d782=8036\", \"49\u00d744=2156\"],\n  [\"65\u00d773=4745\", \"54\u00d779=4266\"],\n  [\"12\u00d781=972\", \"15\u00d723=345\"],\n  [\"13\u00d726=338\", \"64\u00d791=5824\"],\n  [\"82\u00d754=4428\", \"28\u00d789=2492\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication answers throughout the document.\n# Each old string occurs exactly once, so Find/Replace with wdReplaceAll\n# (scoped to the whole document) is safe and unambiguous for every pair.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"19\u00d752=988\", \"33\u00d753=1749\"),\n    @(\"27\u00d727=729\", \"85\u00d788=7480\"),\n    @(\"17\u00d756=952\", \"79\u00d787=6873\"),\n    @(\"44\u00d760=2640\", \"73\u00d718=1314\"),\n    @(\"63\u00d728=1764\", \"21\u00d765=1365\"),\n    @(\"91\u00d754=4914\", \"60\u00d786=5160\"),\n    @(\"85\u00d732=2720\", \"35\u00d789=3115\"),\n    @(\"65\u00d765=4225\", \"68\u00d747=3196\"),\n    @(\"64\u00d720=1280\", \"17\u00d781=1377\"),\n    @(\"32\u00d718=576\", \"34\u00d767=2278\"),\n    @(\"56\u00d743=2408\", \"12\u00d765=780\"),\n    @(\"69\u00d712=828\", \"75\u00d743=3225\"),\n    @(\"39\u00d731=1209\", \"14\u00d757=798\"),\n    @(\"13\u00d781=1053\", \"52\u00d771=3692\"),\n    @(\"39\u00d750=1950\", \"72\u00d713=936\"),\n    @(\"81\u00d758=4698\", \"97\u00d784=8148\"),\n    @(\"36\u00d738=1368\", \"44\u00d766=2904\"),\n    @(\"15\u00d760=900\", \"63\u00d785=5355\"),\n    @(\"39\u00d722=858\", \"82\u00d742=3444\"),\n    @(\"37\u00d757=2109\", \"39\u00d759=2301\"),\n    @(\"98\u00d782=8036\", \"49\u00d744=2156\"),\n    @(\"65\u00d773=4745\", \"54\u00d779=4266\"),\n    @(\"12\u00d781=972\", \"15\u00d723=345\"),\n    @(\"13\u00d726=338\", \"64\u00d791=5824\"),\n    @(\"82\u00d754=4428\", \"28\u00d789=2492\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # wdFindStop = 0, wdReplaceAll = 2\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 0, $false, $newText, 2)\n}\n"}
